# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2304"
#   "<header>_new" -> "<header>_FV2310"
# Then (re)build the header row as an Excel Table ("Table1") covering
# A1:U58, and freeze the header row (pane split below row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = $ws.UsedRange.Columns.Count
$headerRow = 1

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item($headerRow, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2304"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
        }
    }
}

# Turn the header + data range into a native Excel Table, named Table1.
$tableRange = $ws.UsedRange
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1, keep it visible while scrolling).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
